$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H1").Value = "Phone"
$ws.Range("H2").Value = 6285456895136
$ws.Columns.Item(8).ColumnWidth = 11.1
$ws.Range("H11").Select()
